$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the thick-bottom-border spacer row (row 3) to the new columns ---
$ws.Range("J3:K3").Copy()
$ws.Range("L3:M3").PasteSpecial(-4122)

# --- Extend the table with two new year columns (L = 2022, M = 2023) ---
# Copy number/alignment formatting from the existing "2021" column (K) so the
# new columns inherit identical cell styles, then overwrite the handful of
# cells whose style differs from column K.
$ws.Range("K4:K10").Copy()
$ws.Range("L4:L10").PasteSpecial(-4122)
$ws.Range("K4:K10").Copy()
$ws.Range("M4:M10").PasteSpecial(-4122)

# M7 uses the "dash" style (like K6/I6/I8) rather than K7's plain style.
$ws.Range("K6").Copy()
$ws.Range("M7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Header year values ---
$ws.Range("L4").Value = 2022
$ws.Range("M4").Value = 2023

# --- Data rows ---
$ws.Range("L5").Value = 500
$ws.Range("M5").Value = 700

$ws.Range("L6").Value = "-"
$ws.Range("M6").Value = 6

$ws.Range("L7").Value = 23
$ws.Range("M7").Value = "-"

$ws.Range("L8").Value = 7
$ws.Range("M8").Value = 5

$ws.Range("L9").Value = 23
$ws.Range("M9").Value = 21

$ws.Range("L10").Value = 172
$ws.Range("M10").Value = 143

# --- Footnote row (A11:C11) switches to a smaller 8pt Times New Roman font ---
$footnote = $ws.Range("A11:C11")
$footnote.Font.Name = "Times New Roman"
$footnote.Font.Size = 8
$footnote.VerticalAlignment = -4108

# --- Page setup ---
$ws.PageSetup.PaperSize = 256
$ws.PageSetup.Orientation = 1
